$d = $word.ActiveDocument

# The "Type request" table (second table in the document) holds both
# edits: the previously-empty "datetime" description, and the "deadline"
# row that needs to move below "purpose".
$t = $d.Tables.Item(2)

# ---------------------------------------------------------------------
# 1) Fill in the previously-empty Description cell ("Description" column,
#    index 5) on the "datetime" row.
# ---------------------------------------------------------------------
$datetimeRow = $t.Rows.Item(3)
if ($datetimeRow.Cells.Item(1).Range.Text.TrimEnd([char]7,[char]13) -ne "datetime") {
    throw "Expected row 3 to be the 'datetime' row."
}
$datetimeRow.Cells.Item(5).Range.Text = "datetime de création de la demande"

# ---------------------------------------------------------------------
# 2) Re-order the "deadline" row so that it follows "purpose" instead of
#    preceding it. All data rows in this table share the same 6-column
#    layout, so the row is recreated in its new slot: capture its cell
#    text, delete it, insert a fresh row before "freetext", and
#    repopulate that fresh row with the captured text.
# ---------------------------------------------------------------------
$deadlineRow = $t.Rows.Item(5)
if ($deadlineRow.Cells.Item(1).Range.Text.TrimEnd([char]7,[char]13) -ne "deadline") {
    throw "Expected row 5 to be the 'deadline' row."
}
$deadlineCol1 = $deadlineRow.Cells.Item(1).Range.Text
$deadlineCol2 = $deadlineRow.Cells.Item(2).Range.Text
$deadlineCol3 = $deadlineRow.Cells.Item(3).Range.Text
$deadlineCol4 = $deadlineRow.Cells.Item(4).Range.Text
$deadlineCol5 = $deadlineRow.Cells.Item(5).Range.Text

$deadlineRow.Delete()

# After deleting row 5 ("deadline"), the rows shift up: ... convention(4),
# purpose(5), freetext(6). Insert the recreated "deadline" row right
# before "freetext" so it ends up between "purpose" and "freetext".
$freetextRow = $t.Rows.Item(6)
if ($freetextRow.Cells.Item(1).Range.Text.TrimEnd([char]7,[char]13) -ne "freetext") {
    throw "Expected row 6 to be the 'freetext' row after deleting 'deadline'."
}
$newRow = $t.Rows.Add($freetextRow)

$newRow.Cells.Item(1).Range.Text = $deadlineCol1
$newRow.Cells.Item(2).Range.Text = $deadlineCol2
$newRow.Cells.Item(3).Range.Text = $deadlineCol3
$newRow.Cells.Item(4).Range.Text = $deadlineCol4
$newRow.Cells.Item(5).Range.Text = $deadlineCol5
# Cell 6 ("Exemple") was already empty in the source row — leave the
# freshly-inserted blank cell as-is.
